# adds thrust to UI. program now reads dat files when dropped in
#
# Update propeller size labels so they use an "x" separator instead of "*"
# (e.g. "18*10" -> "18x10"), and move the current selection from F29 to
# D24, which also clears the scrolled-down topLeftCell="A4" view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C12").Value = "18x10"
$ws.Range("C13:C23").Value = "19x10"
$ws.Range("C24:C34").Value = "20x10"

$ws.Range("D24").Select()
